# cryptos.xlsx daily refresh
# Updates the "Price" and "Volume(1h)" columns for every coin with new
# market data, and re-syncs rows 47-51 (Coin/Link/Price/Volume) after the
# three lowest-ranked coins shifted position in the source ranking.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel auto-detects plain numeric text (e.g. "291.94") and stores it as a
# number; the sheet keeps the Price column as text, so a leading apostrophe
# (exactly like typing it in the UI) is used to force text storage whenever
# the new value would otherwise be recognized as a number.
function Set-CellText($Range, $Text) {
    if ($Text -match '^-?[0-9]+(\.[0-9]+)?$') {
        $Range.Value = "'" + $Text
    } else {
        $Range.Value = $Text
    }
}

# Row number -> updated column letter/value pairs
$updates = [ordered]@{
    2 = @{ D='43.680.07'; E='  -1.62%  ' }
    3 = @{ D='2.181.17'; E='  -2.78%  ' }
    4 = @{ E='  -0.26%  ' }
    5 = @{ D='291.94'; E='  -4.46%  ' }
    6 = @{ D='86.45'; E='  -6.97%  ' }
    7 = @{ D='0.558'; E='  -2.05%  ' }
    8 = @{ E='  -0.16%  ' }
    9 = @{ D='0.474'; E='  -9.04%  ' }
    10 = @{ D='31.73'; E='  -7.94%  ' }
    11 = @{ D='0.0757'; E='  -6.51%  ' }
    12 = @{ E='  -2.02%  ' }
    13 = @{ D='6.65'; E='  -6.61%  ' }
    14 = @{ D='2.513.62'; E='  -2.86%  ' }
    15 = @{ D='2.249.86'; E='  -4.85%  ' }
    16 = @{ D='12.79'; E='  -5.51%  ' }
    17 = @{ D='0.754'; E='  -9.71%  ' }
    18 = @{ D='43.259.12'; E='  -1.86%  ' }
    19 = @{ D='0.0₃0865'; E='  -9.81%  ' }
    20 = @{ D='5.75'; E='  -9.28%  ' }
    21 = @{ D='10.55'; E='  -14.20%  ' }
    22 = @{ D='62.03'; E='  -5.36%  ' }
    23 = @{ D='226.52'; E='  -4.46%  ' }
    24 = @{ D='2.73'; E='  -6.76%  ' }
    25 = @{ E='  -0.06%  ' }
    26 = @{ D='1.78'; E='  -9.60%  ' }
    27 = @{ E='  +0.46%  ' }
    28 = @{ D='9.05'; E='  -7.53%  ' }
    29 = @{ D='34.49'; E='  -11.08%  ' }
    30 = @{ D='18.85'; E='  -5.75%  ' }
    31 = @{ D='145.02'; E='  -5.71%  ' }
    32 = @{ D='5.18'; E='  -12.55%  ' }
    33 = @{ D='2.47'; E='  -6.99%  ' }
    34 = @{ D='0.0715'; E='  -10.33%  ' }
    35 = @{ D='0.114'; E='  -3.87%  ' }
    36 = @{ D='2.84'; E='  -7.97%  ' }
    37 = @{ D='0.0998'; E='  -7.44%  ' }
    38 = @{ D='1.61'; E='  -11.64%  ' }
    39 = @{ D='1.00'; E='  -0.38%  ' }
    40 = @{ D='0.0275'; E='  -8.31%  ' }
    41 = @{ D='3.43'; E='  -9.73%  ' }
    42 = @{ D='12.94'; E='  -10.75%  ' }
    43 = @{ D='2.98'; E='  -13.03%  ' }
    44 = @{ D='1.749.47'; E='  +1.17%  ' }
    45 = @{ D='1.60'; E='  +1.37%  ' }
    46 = @{ D='71.70'; E='  -10.76%  ' }
    47 = @{ B='Algorand'; C='https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'; D='0.168'; E='  -12.70%  ' }
    48 = @{ B='Aave'; C='https://coinranking.com/coin/ixgUfzmLR+aave-aave'; D='90.26'; E='  -8.99%  ' }
    49 = @{ B='RocketPoolETH'; C='https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'; D='2.396.08'; E='  -2.83%  ' }
    50 = @{ B='HuobiToken'; C='https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'; D='2.69'; E='  +6.92%  ' }
    51 = @{ B='ordi'; C='https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'; D='63.57'; E='  -8.58%  ' }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        Set-CellText $ws.Range("$col$row") $cols[$col]
    }
}

